# Auto-generated Excel COM-interop script to update "想去人数" (F column) values
# across all four worksheets, per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 1205
$ws.Range("F6").Value = 9432
$ws.Range("F7").Value = 7270
$ws.Range("F9").Value = 326
$ws.Range("F10").Value = 5810
$ws.Range("F13").Value = 27
$ws.Range("F14").Value = 6687
$ws.Range("F15").Value = 1119
$ws.Range("F16").Value = 479
$ws.Range("F17").Value = 450
$ws.Range("F19").Value = 656
$ws.Range("F21").Value = 293
$ws.Range("F22").Value = 223
$ws.Range("F25").Value = 10843
$ws.Range("F27").Value = 43
$ws.Range("F28").Value = 2064
$ws.Range("F29").Value = 2616
$ws.Range("F32").Value = 2391
$ws.Range("F33").Value = 92
$ws.Range("F35").Value = 31
$ws.Range("F36").Value = 78
$ws.Range("F38").Value = 1495
$ws.Range("F39").Value = 64
$ws.Range("F40").Value = 25
$ws.Range("F41").Value = 5514
$ws.Range("F43").Value = 777
$ws.Range("F44").Value = 141
$ws.Range("F45").Value = 174
$ws.Range("F47").Value = 1441
$ws.Range("F48").Value = 78
$ws.Range("F49").Value = 1112

$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 50
$ws.Range("F11").Value = 202
$ws.Range("F18").Value = 9
$ws.Range("F20").Value = 35
$ws.Range("F21").Value = 8

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 113
$ws.Range("F3").Value = 204

$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1206
$ws.Range("F5").Value = 9432
$ws.Range("F6").Value = 7270
$ws.Range("F7").Value = 113
$ws.Range("F9").Value = 204
$ws.Range("F12").Value = 5810
$ws.Range("F15").Value = 27
$ws.Range("F16").Value = 6687
$ws.Range("F17").Value = 6687
$ws.Range("F18").Value = 1119
$ws.Range("F19").Value = 479
$ws.Range("F20").Value = 450
$ws.Range("F21").Value = 656
$ws.Range("F23").Value = 293
$ws.Range("F24").Value = 223
$ws.Range("F27").Value = 202
$ws.Range("F28").Value = 10843
$ws.Range("F30").Value = 44
$ws.Range("F31").Value = 2064
$ws.Range("F32").Value = 2616
$ws.Range("F33").Value = 2391
$ws.Range("F34").Value = 92
$ws.Range("F36").Value = 31
$ws.Range("F37").Value = 9
$ws.Range("F39").Value = 1495
$ws.Range("F40").Value = 5514
$ws.Range("F41").Value = 35
$ws.Range("F43").Value = 777
$ws.Range("F44").Value = 141
$ws.Range("F45").Value = 174
$ws.Range("F48").Value = 1441
$ws.Range("F49").Value = 78
$ws.Range("F50").Value = 1112
